$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "29×49=1421" "41×91=3731"
Replace-Text "13×78=1014" "58×50=2900"
Replace-Text "95×37=3515" "28×95=2660"
Replace-Text "17×37=629" "12×12=144"
Replace-Text "24×81=1944" "33×21=693"
Replace-Text "80×75=6000" "31×14=434"
Replace-Text "85×21=1785" "59×69=4071"
Replace-Text "81×37=2997" "41×38=1558"
Replace-Text "75×98=7350" "74×53=3922"
Replace-Text "55×19=1045" "33×59=1947"
Replace-Text "15×82=1230" "28×75=2100"
Replace-Text "21×28=588" "43×71=3053"
Replace-Text "65×24=1560" "90×30=2700"
Replace-Text "43×34=1462" "53×98=5194"
Replace-Text "39×29=1131" "11×35=385"
Replace-Text "95×29=2755" "23×22=506"
Replace-Text "84×69=5796" "20×37=740"
Replace-Text "17×58=986" "82×40=3280"
Replace-Text "92×20=1840" "85×84=7140"
Replace-Text "40×29=1160" "25×74=1850"
Replace-Text "25×71=1775" "63×27=1701"
Replace-Text "77×83=6391" "57×87=4959"
Replace-Text "63×81=5103" "87×20=1740"
Replace-Text "87×22=1914" "72×27=1944"
Replace-Text "50×16=800" "18×56=1008"
